$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1286.2069  # H40: was 1254.6875
$ws.Cells.Item(40, 9).Value = 1268  # I40: was 1225.9259
$ws.Cells.Item(40, 10).Value = 1400  # J40: was 1410
$ws.Cells.Item(40, 11).Value = 1268  # K40: was 1225.9259
$ws.Cells.Item(40, 12).Value = 1400  # L40: was 1410
$ws.Cells.Item(40, 13).Value = -1093  # M40: was -1050.9259
$ws.Cells.Item(40, 14).Value = -1750  # N40: was -1760
$ws.Cells.Item(43, 8).Value = 1561.8182  # H43: was 4811.8213
$ws.Cells.Item(43, 9).Value = 1800  # I43: was 20942.2
$ws.Cells.Item(43, 10).Value = 1538  # J43: was 1305.2174
$ws.Cells.Item(43, 11).Value = 1800  # K43: was 20942.2
$ws.Cells.Item(43, 12).Value = 1538  # L43: was 1305.2174
$ws.Cells.Item(43, 13).Value = -1731  # M43: was -20873.2
$ws.Cells.Item(43, 14).Value = -1676  # N43: was -1443.2174
$ws.Cells.Item(87, 8).Value = 12672.059  # H87: was 13649.782
$ws.Cells.Item(87, 10).Value = 13116.875  # J87: was 14017.728
$ws.Cells.Item(87, 12).Value = 13116.875  # L87: was 14017.728
$ws.Cells.Item(87, 14).Value = -15612.875  # N87: was -16513.728
$ws.Cells.Item(90, 8).Value = 12672.059  # H90: was 13649.782
$ws.Cells.Item(90, 10).Value = 13116.875  # J90: was 14017.728
$ws.Cells.Item(90, 12).Value = 39350.625  # L90: was 42053.18399999999
$ws.Cells.Item(90, 14).Value = -51830.625  # N90: was -54533.18399999999
$ws.Cells.Item(106, 8).Value = 3010.5945  # H106: was 2955.0527
$ws.Cells.Item(106, 9).Value = 2976  # I106: was 2722.5
$ws.Cells.Item(106, 10).Value = 3016  # J106: was 3017.0667
$ws.Cells.Item(106, 11).Value = 2976  # K106: was 2722.5
$ws.Cells.Item(106, 12).Value = 3016  # L106: was 3017.0667
$ws.Cells.Item(106, 13).Value = -2345  # M106: was -2091.5
$ws.Cells.Item(106, 14).Value = -4278  # N106: was -4279.066699999999
$ws.Cells.Item(121, 8).Value = 937.76  # H121: was 1064.9048
$ws.Cells.Item(121, 10).Value = 934.13043  # J121: was 1073.8948
$ws.Cells.Item(121, 12).Value = 2802.39129  # L121: was 3221.6844
$ws.Cells.Item(121, 14).Value = -6296.39129  # N121: was -6715.6844
$ws.Cells.Item(129, 8).Value = 780.72974  # H129: was 790.375
$ws.Cells.Item(129, 10).Value = 855.90625  # J129: was 860.4857
$ws.Cells.Item(129, 12).Value = 2567.71875  # L129: was 2581.4571
$ws.Cells.Item(129, 14).Value = -12567.71875  # N129: was -12581.4571
$ws.Cells.Item(138, 8).Value = 2831.2468  # H138: was 2971.5657
$ws.Cells.Item(138, 9).Value = 1628.826  # I138: was 1756.4048
$ws.Cells.Item(138, 10).Value = 4615.484  # J138: was 4472.647
$ws.Cells.Item(138, 11).Value = 4886.478  # K138: was 5269.2144
$ws.Cells.Item(138, 12).Value = 13846.452  # L138: was 13417.941
$ws.Cells.Item(138, 13).Value = 253.5219999999999  # M138: was -129.2143999999998
$ws.Cells.Item(138, 14).Value = -24126.452  # N138: was -23697.941
$ws.Cells.Item(141, 8).Value = 3925.7144  # H141: was 2172.6667
$ws.Cells.Item(141, 9).Value = 0  # I141: was 2043.5714
$ws.Cells.Item(141, 10).Value = 3925.7144  # J141: was 3980
$ws.Cells.Item(141, 11).Value = 0  # K141: was 6130.7142
$ws.Cells.Item(141, 12).Value = 11777.1432  # L141: was 11940
$ws.Cells.Item(141, 13).ClearContents()  # M141: was -950.7142000000003
$ws.Cells.Item(141, 14).Value = -22137.1432  # N141: was -22300

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5577.49  # H32: was 6755.5
$ws.Cells.Item(32, 9).Value = 4942.297  # I32: was 6040.341
$ws.Cells.Item(32, 11).Value = 4942.297  # K32: was 6040.341
$ws.Cells.Item(32, 13).Value = -4655.297  # M32: was -5753.341
$ws.Cells.Item(61, 8).Value = 2380.02  # H61: was 2258.3584
$ws.Cells.Item(61, 9).Value = 1571.119  # I61: was 1477.762
$ws.Cells.Item(61, 10).Value = 6626.75  # J61: was 5238.8184
$ws.Cells.Item(61, 11).Value = 1571.119  # K61: was 1477.762
$ws.Cells.Item(61, 12).Value = 6626.75  # L61: was 5238.8184
$ws.Cells.Item(61, 13).Value = -1359.119  # M61: was -1265.762
$ws.Cells.Item(61, 14).Value = -7050.75  # N61: was -5662.8184
$ws.Cells.Item(136, 8).Value = 2380.02  # H136: was 2258.3584
$ws.Cells.Item(136, 9).Value = 1571.119  # I136: was 1477.762
$ws.Cells.Item(136, 10).Value = 6626.75  # J136: was 5238.8184
$ws.Cells.Item(136, 11).Value = 4713.357  # K136: was 4433.286
$ws.Cells.Item(136, 12).Value = 19880.25  # L136: was 15716.4552
$ws.Cells.Item(136, 13).Value = -2163.357  # M136: was -1883.286
$ws.Cells.Item(136, 14).Value = -24980.25  # N136: was -20816.4552

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1874.8  # H105: was 1624.3214
$ws.Cells.Item(105, 9).Value = 1882.7273  # I105: was 1603
$ws.Cells.Item(105, 10).Value = 1853  # J105: was 1722.4
$ws.Cells.Item(105, 11).Value = 1882.7273  # K105: was 1603
$ws.Cells.Item(105, 12).Value = 1853  # L105: was 1722.4
$ws.Cells.Item(105, 13).Value = -135.7273  # M105: was 144
$ws.Cells.Item(105, 14).Value = -5347  # N105: was -5216.4

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 100.92308  # H7: was 146.94118
$ws.Cells.Item(7, 9).Value = 85.25  # I7: was 77.818184
$ws.Cells.Item(7, 10).Value = 126  # J7: was 273.66666
$ws.Cells.Item(7, 11).Value = 85.25  # K7: was 77.818184
$ws.Cells.Item(7, 12).Value = 126  # L7: was 273.66666
$ws.Cells.Item(7, 13).Value = 27.75  # M7: was 35.181816
$ws.Cells.Item(7, 14).Value = -352  # N7: was -499.66666
$ws.Cells.Item(31, 8).Value = 4996.048  # H31: was 4454.6
$ws.Cells.Item(31, 9).Value = 6342.1113  # I31: was 5003.9165
$ws.Cells.Item(31, 10).Value = 3986.5  # J31: was 3947.5386
$ws.Cells.Item(31, 11).Value = 6342.1113  # K31: was 5003.9165
$ws.Cells.Item(31, 12).Value = 3986.5  # L31: was 3947.5386
$ws.Cells.Item(31, 13).Value = -6047.1113  # M31: was -4708.9165
$ws.Cells.Item(31, 14).Value = -4576.5  # N31: was -4537.5386
$ws.Cells.Item(34, 8).Value = 4996.048  # H34: was 4454.6
$ws.Cells.Item(34, 9).Value = 6342.1113  # I34: was 5003.9165
$ws.Cells.Item(34, 10).Value = 3986.5  # J34: was 3947.5386
$ws.Cells.Item(34, 11).Value = 6342.1113  # K34: was 5003.9165
$ws.Cells.Item(34, 12).Value = 3986.5  # L34: was 3947.5386
$ws.Cells.Item(34, 13).Value = -6140.1113  # M34: was -4801.9165
$ws.Cells.Item(34, 14).Value = -4390.5  # N34: was -4351.5386

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1179.2759  # H34: was 1243.64
$ws.Cells.Item(34, 10).Value = 1426.95  # J34: was 1589.4375
$ws.Cells.Item(34, 12).Value = 4280.85  # L34: was 4768.3125
$ws.Cells.Item(34, 14).Value = -4448.85  # N34: was -4936.3125
$ws.Cells.Item(39, 8).Value = 1206.2858  # H39: was 1187.8667
$ws.Cells.Item(39, 10).Value = 1208.1666  # J39: was 1186.7693
$ws.Cells.Item(39, 12).Value = 3624.4998  # L39: was 3560.3079
$ws.Cells.Item(39, 14).Value = -4212.4998  # N39: was -4148.3079
$ws.Cells.Item(107, 8).Value = 448.86365  # H107: was 454.79166
$ws.Cells.Item(107, 9).Value = 498.125  # I107: was 473.8889
$ws.Cells.Item(107, 10).Value = 420.7143  # J107: was 443.33334
$ws.Cells.Item(107, 11).Value = 1494.375  # K107: was 1421.6667
$ws.Cells.Item(107, 12).Value = 1262.1429  # L107: was 1330.00002
$ws.Cells.Item(107, 13).Value = 425.625  # M107: was 498.3333
$ws.Cells.Item(107, 14).Value = -5102.1429  # N107: was -5170.000019999999
$ws.Cells.Item(133, 8).Value = 4843.933  # H133: was 5423.933
$ws.Cells.Item(133, 9).Value = 4512.231  # I133: was 4988.25
$ws.Cells.Item(133, 10).Value = 7000  # J133: was 7166.6665
$ws.Cells.Item(133, 11).Value = 13536.693  # K133: was 14964.75
$ws.Cells.Item(133, 12).Value = 21000  # L133: was 21499.9995
$ws.Cells.Item(133, 13).Value = -8476.692999999999  # M133: was -9904.75
$ws.Cells.Item(133, 14).Value = -31120  # N133: was -31619.9995
$ws.Cells.Item(137, 8).Value = 4632585  # H137: was 4905268.5
$ws.Cells.Item(137, 10).Value = 5750250  # J137: was 6176419
$ws.Cells.Item(137, 12).Value = 17250750  # L137: was 18529257
$ws.Cells.Item(137, 14).Value = -17260950  # N137: was -18539457

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(48, 8).Value = 20000  # H48: was 20030
$ws.Cells.Item(48, 10).Value = 20000  # J48: was 20030
$ws.Cells.Item(48, 12).Value = 20000  # L48: was 20030
$ws.Cells.Item(48, 14).Value = -20970  # N48: was -21000
$ws.Cells.Item(69, 8).Value = 28000  # H69: was 0
$ws.Cells.Item(69, 10).Value = 28000  # J69: was 0
$ws.Cells.Item(69, 12).Value = 28000  # L69: was 0
$ws.Cells.Item(69, 14).Value = -29498  # N69: was None
$ws.Cells.Item(72, 8).Value = 28000  # H72: was 0
$ws.Cells.Item(72, 10).Value = 28000  # J72: was 0
$ws.Cells.Item(72, 12).Value = 84000  # L72: was 0
$ws.Cells.Item(72, 14).Value = -91488  # N72: was None

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 621.0714  # H16: was 940.1429000000001
$ws.Cells.Item(16, 9).Value = 621.0714  # I16: was 940.1429000000001
$ws.Cells.Item(16, 11).Value = 621.0714  # K16: was 940.1429000000001
$ws.Cells.Item(16, 13).Value = -451.0714  # M16: was -770.1429000000001
$ws.Cells.Item(46, 8).Value = 157851.28  # H46: was 209636.38
$ws.Cells.Item(46, 9).Value = 1244.3334  # I46: was 3750
$ws.Cells.Item(46, 10).Value = 245942.69  # J46: was 264539.4
$ws.Cells.Item(46, 11).Value = 1244.3334  # K46: was 3750
$ws.Cells.Item(46, 12).Value = 245942.69  # L46: was 264539.4
$ws.Cells.Item(46, 13).Value = -1056.3334  # M46: was -3562
$ws.Cells.Item(46, 14).Value = -246318.69  # N46: was -264915.4
$ws.Cells.Item(68, 8).Value = 2945.611  # H68: was 2790.3416
$ws.Cells.Item(68, 9).Value = 1242  # I68: was 1301
$ws.Cells.Item(68, 10).Value = 2994.2856  # J68: was 2951.3513
$ws.Cells.Item(68, 11).Value = 1242  # K68: was 1301
$ws.Cells.Item(68, 12).Value = 2994.2856  # L68: was 2951.3513
$ws.Cells.Item(68, 13).Value = -493  # M68: was -552
$ws.Cells.Item(68, 14).Value = -4492.2856  # N68: was -4449.3513
$ws.Cells.Item(71, 8).Value = 2945.611  # H71: was 2790.3416
$ws.Cells.Item(71, 9).Value = 1242  # I71: was 1301
$ws.Cells.Item(71, 10).Value = 2994.2856  # J71: was 2951.3513
$ws.Cells.Item(71, 11).Value = 6210  # K71: was 6505
$ws.Cells.Item(71, 12).Value = 14971.428  # L71: was 14756.7565
$ws.Cells.Item(71, 13).Value = -2466  # M71: was -2761
$ws.Cells.Item(71, 14).Value = -22459.428  # N71: was -22244.7565
$ws.Cells.Item(93, 8).Value = 1667.7333  # H93: was 3168
$ws.Cells.Item(93, 9).Value = 1232  # I93: was 0
$ws.Cells.Item(93, 10).Value = 4500  # J93: was 3168
$ws.Cells.Item(93, 11).Value = 1232  # K93: was 0
$ws.Cells.Item(93, 12).Value = 4500  # L93: was 3168
$ws.Cells.Item(93, 13).Value = 16  # M93: was None
$ws.Cells.Item(93, 14).Value = -6996  # N93: was -5664
$ws.Cells.Item(136, 8).Value = 2679.7017  # H136: was 2816.4
$ws.Cells.Item(136, 9).Value = 2118.575  # I136: was 2240.8647
$ws.Cells.Item(136, 10).Value = 4000  # J136: was 3999.4443
$ws.Cells.Item(136, 11).Value = 6355.724999999999  # K136: was 6722.5941
$ws.Cells.Item(136, 12).Value = 12000  # L136: was 11998.3329
$ws.Cells.Item(136, 13).Value = -3805.724999999999  # M136: was -4172.5941
$ws.Cells.Item(136, 14).Value = -17100  # N136: was -17098.3329

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(38, 8).Value = 30000  # H38: was 27764
$ws.Cells.Item(38, 9).Value = 0  # I38: was 1056
$ws.Cells.Item(38, 10).Value = 30000  # J38: was 36666.668
$ws.Cells.Item(38, 11).Value = 0  # K38: was 1056
$ws.Cells.Item(38, 12).Value = 30000  # L38: was 36666.668
$ws.Cells.Item(38, 13).ClearContents()  # M38: was -583
$ws.Cells.Item(38, 14).Value = -30946  # N38: was -37612.668
$ws.Cells.Item(132, 8).Value = 812.1688  # H132: was 877.13434
$ws.Cells.Item(132, 9).Value = 723.9355  # I132: was 771.7778
$ws.Cells.Item(132, 10).Value = 1176.8667  # J132: was 1314.7693
$ws.Cells.Item(132, 11).Value = 2171.8065  # K132: was 2315.3334
$ws.Cells.Item(132, 12).Value = 3530.6001  # L132: was 3944.3079
$ws.Cells.Item(132, 13).Value = 358.1934999999999  # M132: was 214.6666
$ws.Cells.Item(132, 14).Value = -8590.6001  # N132: was -9004.3079
$ws.Cells.Item(136, 8).Value = 818.3226  # H136: was 736.05
$ws.Cells.Item(136, 9).Value = 772.9091  # I136: was 683.1613
$ws.Cells.Item(136, 10).Value = 929.3333  # J136: was 918.2222
$ws.Cells.Item(136, 11).Value = 2318.7273  # K136: was 2049.4839
$ws.Cells.Item(136, 12).Value = 2787.9999  # L136: was 2754.6666
$ws.Cells.Item(136, 13).Value = 231.2727  # M136: was 500.5160999999998
$ws.Cells.Item(136, 14).Value = -7887.9999  # N136: was -7854.6666
